$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New course row
$ws.Range("A14").Value = "Writing Efficient Python Code"
$ws.Range("B14").Value = 5

# Copy the font formatting from A13 (dark blue "new" courses) onto A14
$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial(-4122)  # xlPasteFormats

Write-Host "done"
